$d = $word.ActiveDocument

$replacements = @(
    @{old="66×42="; new="66×89="},
    @{old="87×97="; new="77×72="},
    @{old="84×34="; new="48×73="},
    @{old="26×43="; new="67×80="},
    @{old="59×33="; new="31×47="},
    @{old="27×96="; new="41×24="},
    @{old="95×31="; new="51×65="},
    @{old="90×11="; new="60×64="},
    @{old="88×13="; new="64×81="},
    @{old="68×41="; new="50×84="},
    @{old="81×62="; new="94×37="},
    @{old="80×84="; new="19×65="},
    @{old="82×99="; new="36×61="},
    @{old="94×55="; new="56×28="},
    @{old="75×64="; new="74×96="},
    @{old="16×78="; new="21×96="},
    @{old="85×12="; new="26×63="},
    @{old="64×97="; new="69×34="},
    @{old="93×47="; new="56×67="},
    @{old="89×44="; new="46×14="},
    @{old="74×94="; new="11×54="},
    @{old="32×90="; new="69×15="},
    @{old="33×88="; new="81×33="},
    @{old="55×26="; new="41×13="},
    @{old="36×82="; new="52×55="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
